$q4Data = @(
    @(0, "159883", "永赢中证全指医疗器械ETF", "22.83", "99.25", "2.67", "0.6096", 6),
    @(1, "512100", "南方中证1000ETF", "106.09", "98.15", "0.31", "0.3289", 9),
    @(2, "560010", "广发中证1000ETF", "67.21", "98.32", "0.32", "0.2151", 9),
    @(3, "159845", "华夏中证1000ETF", "62.11", "98.70", "0.32", "0.1988", 9),
    @(4, "159629", "富国中证1000ETF", "59.33", "99.34", "0.32", "0.1899", 9),
    @(5, "159633", "易方达中证1000ETF", "58.64", "98.77", "0.32", "0.1876", 9),
    @(6, "515860", "嘉实新兴科技100ETF", "2.04", "98.57", "3.77", "0.0769", 7),
    @(7, "159873", "天弘中证全指医疗保健设备与服务ETF", "2.08", "99.12", "2.25", "0.0468", 8),
    @(8, "159898", "招商中证全指医疗器械ETF", "1.70", "99.16", "2.69", "0.0457", 6),
    @(9, "516790", "华泰柏瑞中证全指医疗保健ETF", "1.61", "98.06", "2.25", "0.0362", 8),
    @(10, "159891", "建信中证全指医疗保健设备与服务ETF", "1.40", "98.79", "2.25", "0.0315", 8),
    @(11, "159877", "南方中证全指医疗保健ETF", "0.98", "99.57", "2.24", "0.0220", 8),
    @(12, "516610", "大成中证全指医疗保健设备与服务ETF", "0.78", "97.46", "2.23", "0.0174", 8),
    @(13, "560110", "汇添富中证1000ETF", "5.61", "93.96", "0.30", "0.0168", 9),
    @(14, "159797", "汇添富中证全指医疗器械ETF", "0.56", "98.45", "2.65", "0.0148", 6),
    @(15, "516300", "华泰柏瑞中证1000ETF", "1.57", "98.94", "0.32", "0.0050", 9),
    @(16, "501069", "华宝标普中国Ａ股质量价值指数（LOF）", "0.14", "93.83", "2.46", "0.0034", 10),
    @(17, "162413", "华宝中证1000指数A", "0.40", "91.60", "0.29", "0.0012", 9),
    @(18, "016033", "华宝中证1000指数C", "0.09", "91.60", "0.29", "0.0003", 9),
)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by duplicating "2022-Q3" (so it inherits
#    identical formatting), then position it right before "2022-Q3".
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

# The duplicated sheet has 6 rows (1 header + 5 data rows); we need 20 rows
# (1 header + 19 data rows), so insert 14 additional rows before row 7.
$newSheet.Rows("7:20").Insert()
$newSheet.Range("A7:A20").ClearFormats()
$newSheet.Range("B7:H20").ClearFormats()
$newSheet.Range("A2").Copy()
$newSheet.Range("A7:A20").PasteSpecial(-4122)

$r = 2
foreach ($row in $q4Data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
# Drop the quote-prefix text-number style so the text cells carry no extra
# number format (matches plain inlineStr cells with no style in the source).
$newSheet.Range("B2:G20").ClearFormats()

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q4 above
#    the existing 2022-Q3 entry and renumber the index column.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()
$total.Rows(2).ClearFormats()

# Restore the index-column style (bold/centered/bordered) on the new A2 cell
# by copying the format from the (now shifted) A3 cell, which already has it.
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 19
$total.Cells.Item(2, 4).Value = 2.05

# Renumber the sequential index in column A for the rows that shifted down.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3

Write-Host "Edit complete"
